# Daily attendance processing - 2025-12-17 11:00:43
# Normalizes the "Recorded By" (column G) values so that the "System"
# token, when it is the first entry in the comma-separated list, is
# moved to the end of the list instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ",\s*"

        if ($parts.Count -gt 1 -and $parts[0] -eq "System") {
            $rest = $parts[1..($parts.Count - 1)]
            $newParts = $rest + @("System")
            $newValue = [string]::Join(", ", $newParts)
            $cell.Value2 = $newValue
        }
    }
}
